$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2881169905109251
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 9506434351236724736.0
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 9506434351251499008.0

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 20.15985084044064

$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 9.983522426115931
$ws.Range("D4").Value = 18.71679738969934
$ws.Range("E4").Value = 13540782231280939008.0
$ws.Range("G4").Value = 13540782231280939008.0

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 9.983522426115931
$ws.Range("D5").Value = 3.223369029078222
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 30.34306516417429
